# "almost finish transport order flow"
#
# Extends the transport-order test-data sheet with 15 new columns
# (Q:AE) covering goods/tunnel/contact/payment details, fills in their
# row-2 sample values, and tweaks two existing row-2 values
# (ORDER_END_DATE -> "Tomorrow", VEHICLE_TYPE -> "Van").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a cell so it ends up as a plain shared-string (t="s"),
# never auto-coerced to a boolean/number, by entering it as a quoted
# text formula and then flattening the formula down to its value via
# copy / paste-values. PasteSpecial(formats) first (copied from a
# same-row "template" cell) carries over the exact existing style, so
# no new style entries are added to styles.xml.
function Set-TextCell($addr, $templateAddr, $text) {
    $ws.Range($templateAddr).Copy()
    $ws.Range($addr).PasteSpecial(-4122)   # xlPasteFormats

    $escaped = $text.Replace('"', '""')
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)   # xlPasteValues
}

# ---------------------------------------------------------------------
# 1) New header row (row 1), columns Q..AE - same style as P1.
# ---------------------------------------------------------------------
$headers = @(
    @("Q1", "GOODS_LONGER_THAN_6_FT"),
    @("R1", "GOODS_TALLER_THAN_2_FT"),
    @("S1", "PET_FRIENDLY"),
    @("T1", "ENGLISH_SPEAKING"),
    @("U1", "TUNNEL_PREFERENCE"),
    @("V1", "SELECTED_TUNNEL"),
    @("W1", "MOVE_DOOR_TO_DOOR"),
    @("X1", "TRANSPORT_OR_DISPOSE_WASTE"),
    @("Y1", "USER_NAME"),
    @("Z1", "PHONE_NUMBER"),
    @("AA1", "EXTENSION_FLAG"),
    @("AB1", "EXTENSION"),
    @("AC1", "TIP"),
    @("AD1", "COUPON"),
    @("AE1", "PAYMENT_METHOD")
)

foreach ($pair in $headers) {
    Set-TextCell $pair[0] "P1" $pair[1]
}

# ---------------------------------------------------------------------
# 2) New row-2 (data row) cells, columns Q..AE - same style as P2.
# ---------------------------------------------------------------------
$data = @(
    @("Q2", "TRUE"),
    @("R2", "TRUE"),
    @("S2", "TRUE"),
    @("T2", "TRUE"),
    @("U2", "TRUE"),
    @("V2", "Tai Lam Tunnel"),
    @("W2", "TRUE"),
    @("X2", "TRUE"),
    @("Y2", "William Koh"),
    @("Z2", "51112222"),
    @("AA2", "TRUE"),
    @("AB2", "852"),
    @("AC2", "20"),
    @("AD2", "FALSE"),
    @("AE2", "FPS")
)

foreach ($pair in $data) {
    Set-TextCell $pair[0] "P2" $pair[1]
}

# ---------------------------------------------------------------------
# 3) Tweak two existing row-2 values (plain text, same cell/style).
# ---------------------------------------------------------------------
Set-TextCell "J2" "J2" "Tomorrow"
Set-TextCell "M2" "M2" "Van"

$excel.CutCopyMode = 0
